$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly cryptos-list refresh: update the Price (D) and Volume(1h) (E)
# columns for each coin row. Price strings such as "1.002" or "44.80"
# must stay text (matching the original inline-string cells), so force
# the Text number format before writing them -- otherwise Excel's
# smart-parse would silently convert them to numbers and drop the
# trailing zero / formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.243.91"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.812.21"
$ws.Range("E3").Value = "  +4.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.29"
$ws.Range("E5").Value = "  +2.23%  "

$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4455"
$ws.Range("E7").Value = "  +5.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3702"
$ws.Range("E8").Value = "  +3.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.80"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07704"
$ws.Range("E10").Value = "  +4.38%  "

$ws.Range("E11").Value = "  +2.22%  "

$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("E13").Value = "  +3.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.299"
$ws.Range("E14").Value = "  +4.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.569"
$ws.Range("E15").Value = "  +6.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.843.72"
$ws.Range("E16").Value = "  +6.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.03"
$ws.Range("E17").Value = "  +7.44%  "

$ws.Range("E18").Value = "  +2.86%  "

$ws.Range("E19").Value = "  +10.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("E21").Value = "  +5.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.228"
$ws.Range("E22").Value = "  +3.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.296.50"
$ws.Range("E23").Value = "  +3.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.70"
$ws.Range("E24").Value = "  +3.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.144"
$ws.Range("E25").Value = "  -10.35%  "

$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.44"
$ws.Range("E27").Value = "  +5.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.033.60"
$ws.Range("E28").Value = "  +5.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.320"
$ws.Range("E29").Value = "  -0.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.48"
$ws.Range("E30").Value = "  +2.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.200"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.895"
$ws.Range("E32").Value = "  +5.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09263"
$ws.Range("E33").Value = "  +2.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.658"
$ws.Range("E34").Value = "  +1.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.10"
$ws.Range("E35").Value = "  +3.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02356"
$ws.Range("E36").Value = "  +5.19%  "

$ws.Range("E37").Value = "  +1.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.182"
$ws.Range("E38").Value = "  +3.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06238"
$ws.Range("E39").Value = "  +2.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6580"
$ws.Range("E40").Value = "  +4.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.197"
$ws.Range("E41").Value = "  +1.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.150"
$ws.Range("E42").Value = "  +3.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.400"
$ws.Range("E44").Value = "  -0.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.92"
$ws.Range("E45").Value = "  +3.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6090"
$ws.Range("E46").Value = "  +5.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.772"
$ws.Range("E47").Value = "  +1.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.98"
$ws.Range("E48").Value = "  +2.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.037"
$ws.Range("E49").Value = "  +5.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.159"
$ws.Range("E50").Value = "  +5.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06987"
$ws.Range("E51").Value = "  +2.67%  "
